$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates
$ws.Range("D1").Value = "iOS App Rating"
$ws.Range("E1").Value = "iOS Review Count"
$ws.Range("F1").Value = "iOS Rank"

# Update date/time column B for rows 2..38 to the new timestamp
$newDate = "2023-06-26 18:41:30"
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 2).Value = $newDate
}

# Update Review Count (E) and Rank (F) values per-row
$ws.Range("E2").Value = 61742
$ws.Range("E3").Value = 4426
$ws.Range("E5").Value = 399031
$ws.Range("F5").Value = 126
$ws.Range("E6").Value = 4754388
$ws.Range("F6").Value = 9
$ws.Range("E7").Value = 43149
$ws.Range("E8").Value = 322237
$ws.Range("E9").Value = 2412812
$ws.Range("E10").Value = 42796
$ws.Range("F10").Value = 94
$ws.Range("E11").Value = 24958
$ws.Range("F11").Value = 51
$ws.Range("E12").Value = 2010064
$ws.Range("F12").Value = 12
$ws.Range("E13").Value = 1047309
$ws.Range("E14").Value = 258060
$ws.Range("F14").Value = 118
$ws.Range("E15").Value = 116043
$ws.Range("E16").Value = 77035
$ws.Range("F16").Value = 165
$ws.Range("E18").Value = 411629
$ws.Range("E19").Value = 3991303
$ws.Range("F19").Value = 14
$ws.Range("E21").Value = 6603
$ws.Range("F22").Value = 58
$ws.Range("E23").Value = 45221
$ws.Range("E24").Value = 876975
$ws.Range("F24").Value = 7
$ws.Range("E25").Value = 291

# Row 28: new company entry (Astound Mobile) replacing the former "App Store" placeholder row
$ws.Range("C28").Value = "Astound Mobile"
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 1

$ws.Range("E33").Value = 28371
$ws.Range("E35").Value = 13437
$ws.Range("E37").Value = 3298
